$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure columns -------------------------------------------------
# Current layout: A=evaluation, B=description, C=last_run
# Target layout:  A=id_evaluation (new), B=evaluation (was A), C=description (was B, repurposed),
#                 D=output (new), E=last_run (was C)
$ws.Columns("A").Insert()
$ws.Columns("D").Insert()

# --- Header row -----------------------------------------------------------
$ws.Cells.Item(1,1).Value = "id_evaluation"
$ws.Cells.Item(1,2).Value = "evaluation"
$ws.Cells.Item(1,3).Value = "description"
$ws.Cells.Item(1,4).Value = "output"
$ws.Cells.Item(1,5).Value = "last_run"

# --- Data rows --------------------------------------------------------------
$names = @("SIMCE Matemáticas","SIMCE Lenguaje","DIA Matemáticas","DIA Lenguaje","Cálculo Veloz","Fluidez Lectora","En Pullinque Todos Leemos","PDL","DIA Extraer respuestas correctas")
$descriptions = @("Workflow SIMCE","Workflow SIMCE","Workflow DIA","","","","","","")
$outputs = @("PDF","PDF","PDF","PDF","PDF","PDF","PDF","PDF","Excel")

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r,1).Value = $i + 1
    $ws.Cells.Item($r,2).Value = $names[$i]
    if ($descriptions[$i] -ne "") {
        $ws.Cells.Item($r,3).Value = $descriptions[$i]
    }
    $ws.Cells.Item($r,4).Value = $outputs[$i]
    $ws.Cells.Item($r,5).NumberFormat = "dd\ mmm\ yyyy"
    $ws.Cells.Item($r,5).Value = 46050
}

# --- Column widths ----------------------------------------------------------
$ws.Columns("B").ColumnWidth = 30.083333333333332
$ws.Columns("C").ColumnWidth = 32.583333333333336
$ws.Columns("D").ColumnWidth = 18.083333333333332
$ws.Columns("E").ColumnWidth = 10.416666666666666

# --- Selection --------------------------------------------------------------
$ws.Range("A11").Select() | Out-Null
